$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Đơn phụ phẫu 1" (2nd sheet): a new service-order row is inserted as
# row 4 (shifting the previous "Tổng" summary row down to row 5) and the
# summary row's totals are updated to account for the new entry.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Push the current row 4 (old "Tổng" summary row) down to row 5, preserving
# its existing (empty) formatting for the untouched columns C:H.
$ws2.Range("A4:I4").Copy($ws2.Range("A5:I5"))

# Update the summary totals on the (now relocated) row 5.
$ws2.Range("B5").Value = 3
$ws2.Range("I5").Value = 100000

# Overwrite row 4 with the new service-order data.
$ws2.Range("A4").Value = "HD-LUXURY"
$ws2.Range("B4").Value = 680
$ws2.Range("C4").NumberFormat = "@"
$ws2.Range("C4").Value = "08-24-2024"
$ws2.Range("D4").Value = "SÓC TRĂNG"
$ws2.Range("E4").Value = "triệu tú kiều "
$ws2.Range("F4").Value = "Cá nhân"
$ws2.Range("G4").Value = "Cắt mí"
$ws2.Range("H4").Value = "Trần Khánh Hiệp"
$ws2.Range("I4").Value = 50000

# ---------------------------------------------------------------------------
# Sheet "Lương" (4th sheet): reflect the extra 50 000 of "Công phụ phẫu 1"
# pay at SÓC TRĂNG and the resulting change to the total salary figures.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("B29").Value = 100000
$ws4.Range("B34").Value = -3671428.571428571
$ws4.Range("B35").Value = -3671428.571428571
